# Insert a new data row above row 162 (shifts the existing rows 162-184
# down to 163-185) and populate it with the new Berenjena record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 162:184 down by one row, preserving formatting/styles.
$ws.Range("A162").EntireRow.Insert()

# Fill in the freshly inserted row 162 with its data.
$ws.Range("A162").Value = 5
$ws.Range("B162").Value = "Macroferia Regional de Talca"
$ws.Range("C162").Value = "Maule"
$ws.Range("D162").Value = 45127
$ws.Range("E162").Value = 7
$ws.Range("F162").Value = 100112001
$ws.Range("G162").Value = "Berenjena"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 200
$ws.Range("K162").Value = 7000
$ws.Range("L162").Value = 7000
$ws.Range("M162").Value = 7000
$ws.Range("N162").Value = "`$/caja 50 unidades"
$ws.Range("O162").Value = "Región de Arica y Parinacota"
$ws.Range("P162").Value = 140
$ws.Range("Q162").Value = 50
$ws.Range("R162").Value = "Hortaliza"
